$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "q" parameter row marker at A20 (per commit "Add uniform disc source")
$ws.Range("A20").Value = "q"

# Change D5:D10 Type from "Parameterised TNSA" to "Gaussian" (uniform disc -> Gaussian source)
$ws.Range("D5:D10").Value = "Gaussian"

$ws.Range("D10").Select()
